$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Rita", "rita@gmail.com", 1),
    @("Rita", "rita@gmail.com", 1),
    @("sasha", "sasha@gmail.com", 9),
    @("dasha", "dasha@gnail.com", 0),
    @("Nadja", "nadja@gmail.com", 0)
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 3)).Style = "Normal"
    $row++
}
